$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2023-09-10 -> 2023-09-11, i.e. serial 45179 -> 45180) for every data row
# (rows 2 through 410).
$ws.Range("C2:C410").Value = 45180
